# Updated cryptos list (Price / Volume(1h) columns) to match latest scrape.
# Numeric-looking Price values are written with a leading apostrophe so
# Excel stores them as literal text (matching the source data, which keeps
# values such as "211.45" or "19.47" as plain strings) instead of silently
# re-parsing them as floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.643.30"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "1.596.39"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'211.45"

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.246"
$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("D10").Value = "'19.47"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "1.821.97"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "1.587.27"
$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Value = "26.643.91"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").Value = "'208.74"
$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("E21").Value = "  +4.13%  "

$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("E23").Value = "  +2.94%  "

$ws.Range("D24").Value = "'8.99"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").Value = "'143.86"
$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").Value = "'15.28"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "'0.0513"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  +0.81%  "

$ws.Range("D34").Value = "1.289.48"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").Value = "'0.618"
$ws.Range("E35").Value = "  -7.00%  "

$ws.Range("E36").Value = "  +0.49%  "

$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D39").Value = "'0.830"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("D40").Value = "'1.02"
$ws.Range("E40").Value = "  +14.35%  "

$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("D43").Value = "'0.782"
$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("D44").Value = "'63.21"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("D45").Value = "1.733.16"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").Value = "'91.08"
$ws.Range("E46").Value = "  +1.03%  "

$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "  -2.92%  "

$ws.Range("E48").Value = "  +0.89%  "

$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").Value = "'7.33"
$ws.Range("E51").Value = "  -1.81%  "
